$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: new task entry ---
# A11: date text "04.07.2019" stored as plain text (not a date serial number)
$ws.Range("A11").Value = "'04.07.2019"
$ws.Range("A11").Style = "Normal"

# B11: description text, wrapped like the other description cells
$ws.Range("B11").Value = "Added fonts for displaying various Mtg symbols such as`nmanacosts and card types"
$ws.Range("B11").WrapText = $true

# C11: numeric hours value
$ws.Range("C11").Value = 1

# Row height matches the other wrapped, two-line rows (30pt)
$ws.Rows.Item(11).RowHeight = 30

# Update selection to reflect the newly active cell
$ws.Range("B11").Select()
